# extraction des eleves ayant plus de 20ans
# Convert the "Age" column (E) from text like "18 ans" to a plain number (18)
# for every student row, so the age can be used for numeric filtering/extraction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 59

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $text = $cell.Value()
    if ($text -ne $null -and $text -ne "") {
        $number = [int]($text.ToString().Split(" ")[0])
        $cell.Value = $number
    }
}

# This row's first-name cell was cleared while the age value was being entered
$ws.Range("B33").Value = ""

# Leave the selection where the last edit happened
$ws.Range("E28").Select()
